$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57:74 down to 58:75,
# preserving their values/formats (matches the diff's net effect of a new
# weekly record being inserted into the dataset).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record. The
# non-varying descriptive columns (A,B,C,E,F,G,H,I,J,K,L,T) mirror the
# constant values used throughout this block of rows.
$ws.Cells.Item(57, 1).Value = 9
$ws.Cells.Item(57, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(57, 3).Value = "Metropolitana"
$ws.Cells.Item(57, 4).Value = 44754
$ws.Cells.Item(57, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(57, 5).Value = 13
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100102
$ws.Cells.Item(57, 8).Value = "Cítricos"
$ws.Cells.Item(57, 9).Value = 100102006
$ws.Cells.Item(57, 10).Value = "Pomelo"
$ws.Cells.Item(57, 11).Value = "Start Ruby"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 400
$ws.Cells.Item(57, 14).Value = 8500
$ws.Cells.Item(57, 15).Value = 8500
$ws.Cells.Item(57, 16).Value = 8500
$ws.Cells.Item(57, 17).Value = "$/caja 14 kilos"
$ws.Cells.Item(57, 18).Value = "Región Metropolitana"
$ws.Cells.Item(57, 19).Value = 607
$ws.Cells.Item(57, 20).Value = 14
